$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells B2:C6 with new ticker symbols
$ws.Range("B2").Value = "NSE:COLPAL"
$ws.Range("C2").Value = "NSE:AXISILVER"

$ws.Range("B3").Value = "NSE:DIVISLAB"
$ws.Range("C3").Value = "NSE:INDOTHAI"

$ws.Range("B4").Value = "NSE:DMART"
$ws.Range("C4").Value = "NSE:MEDPLUS"

$ws.Range("B5").Value = "NSE:GODREJCP"
$ws.Range("C5").Value = "NSE:RVHL"

$ws.Range("B6").Value = "NSE:HINDUNILVR"

# Add new rows 7-9
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "NSE:JUBLFOOD"

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "NSE:PGHH"

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "NSE:PGHL"

# Apply style from A2:A6 (style index 1) to the new A7:A9 cells
$ws.Range("A2").Copy()
$ws.Range("A7:A9").PasteSpecial(-4122)

# Materialize the trailing empty-text cells (C:F) on rows 7-9, mirroring
# the blank inline-string cells that already exist on rows 2-6. A bare
# Value = "" is treated as "no-op" by the engine (it won't create a cell
# for a fully blank write), so force text-typed blank content via a
# quote-prefixed empty literal and then reset the style back to Normal
# so no stray "quote prefix" formatting is left behind.
foreach ($r in 7..9) {
    foreach ($col in @("C", "D", "E", "F")) {
        $cell = $ws.Range("$col$r")
        $cell.Value = "'"
        $cell.Style = "Normal"
    }
}
